$d = $word.ActiveDocument
Write-Host ("CustomXMLParts.Count = " + $d.CustomXMLParts.Count)
for ($i = 1; $i -le $d.CustomXMLParts.Count; $i++) {
    $part = $d.CustomXMLParts.Item($i)
    Write-Host ("Part " + $i + ": " + $part.NamespaceURI)
}
